$d = $word.ActiveDocument

# --- Step 1: turn " 1" into " " + "2" as separate runs -------------------
# Using TrackRevisions while editing makes Word split the run boundaries
# the same way a human edit would (leaving "First ", "demo", " " intact and
# isolating the replaced character into its own run), then AcceptAllRevisions
# bakes the edit in without leaving any <w:ins>/<w:del> markup behind.
$d.TrackRevisions = $true
$d.Range(11, 12).Text = "2"
$d.AcceptAllRevisions()
$d.TrackRevisions = $false

# --- Step 2: add the new paragraph "Thêm một dòng" ------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs(2)
$p2.Range.InsertAfter("Thêm một dòng")

# --- Step 3: relocate the _GoBack bookmark into the new last paragraph ---
# It needs to end up collapsed right after the new text, before the final
# paragraph mark. Adding a bookmark collapsed exactly at the very end of the
# document content needs a trailing character present, so a throwaway
# character is appended first, the bookmark is planted just before it, and
# then the throwaway character is removed again.
$d.Bookmarks("_GoBack").Delete()

$endPos = $d.Content.End
$d.Range($endPos - 1, $endPos - 1).InsertAfter("X")

$targetPos = $endPos - 1
$bmRange = $d.Range($targetPos, $targetPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

$endPos2 = $d.Content.End
$d.Range($endPos2 - 2, $endPos2 - 1).Delete()
